$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1673151750972763
$ws.Range("C2").Value = 0.5836575875486382
$ws.Range("J2").Value = 0.02723735408560311
$ws.Range("P2").Value = 0.1284046692607004
$ws.Range("S2").Value = 0.0933852140077821
$ws.Range("C3").Value = 0.03184713375796178
$ws.Range("J3").Value = 0.05095541401273886
$ws.Range("P3").Value = 0.732484076433121
$ws.Range("S3").Value = 0.1847133757961783
$ws.Range("J4").Value = 0.06122448979591837
$ws.Range("P4").Value = 0.6938775510204082
$ws.Range("S4").Value = 0.2448979591836735
$ws.Range("B6").Value = 0.04633204633204633
$ws.Range("D6").Value = 0.007722007722007722
$ws.Range("E6").Value = 0.003861003861003861
$ws.Range("F6").Value = 0.04247104247104247
$ws.Range("J6").Value = 0.1930501930501931
$ws.Range("O6").Value = 0.01544401544401544
$ws.Range("Q6").Value = 0.1505791505791506
$ws.Range("R6").Value = 0.1081081081081081
$ws.Range("S6").Value = 0.4324324324324325
$ws.Range("B7").Value = 0.1055045871559633
$ws.Range("D7").Value = 0.01834862385321101
$ws.Range("E7").Value = 0.009174311926605505
$ws.Range("F7").Value = 0.06880733944954129
$ws.Range("J7").Value = 0.09174311926605505
$ws.Range("O7").Value = 0.02293577981651376
$ws.Range("Q7").Value = 0.1559633027522936
$ws.Range("R7").Value = 0.0871559633027523
$ws.Range("S7").Value = 0.4403669724770642
$ws.Range("B8").Value = 0.096
$ws.Range("D8").Value = 0.02
$ws.Range("E8").Value = 0.002
$ws.Range("F8").Value = 0.07000000000000001
$ws.Range("J8").Value = 0.102
$ws.Range("O8").Value = 0.024
$ws.Range("Q8").Value = 0.166
$ws.Range("R8").Value = 0.12
$ws.Range("S8").Value = 0.4
$ws.Range("B9").Value = 0.08
$ws.Range("D9").Value = 0.02285714285714286
$ws.Range("F9").Value = 0.06857142857142857
$ws.Range("J9").Value = 0.1085714285714286
$ws.Range("O9").Value = 0.01142857142857143
$ws.Range("Q9").Value = 0.1885714285714286
$ws.Range("R9").Value = 0.1142857142857143
$ws.Range("S9").Value = 0.4057142857142857
$ws.Range("B10").Value = 0.09010458567980692
$ws.Range("D10").Value = 0.02333065164923572
$ws.Range("F10").Value = 0.09090909090909091
$ws.Range("J10").Value = 0.1005631536604988
$ws.Range("O10").Value = 0.02413515687851971
$ws.Range("Q10").Value = 0.1810136765888978
$ws.Range("R10").Value = 0.1134352373290426
$ws.Range("S10").Value = 0.3765084473049075
$ws.Range("G11").Value = 0.1219512195121951
$ws.Range("J11").Value = 0.09451219512195122
$ws.Range("K11").Value = 0.1829268292682927
$ws.Range("L11").Value = 0.5823170731707317
$ws.Range("S11").Value = 0.01829268292682927
$ws.Range("G12").Value = 0.7564766839378239
$ws.Range("J12").Value = 0.1968911917098446
$ws.Range("K12").Value = 0.005181347150259068
$ws.Range("L12").Value = 0.01036269430051814
$ws.Range("S12").Value = 0.0310880829015544
$ws.Range("G13").Value = 0.7115384615384616
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.03846153846153846
$ws.Range("F15").Value = 0.01149425287356322
$ws.Range("H15").Value = 0.1494252873563219
$ws.Range("I15").Value = 0.06130268199233716
$ws.Range("J15").Value = 0.3371647509578544
$ws.Range("K15").Value = 0.05747126436781609
$ws.Range("M15").Value = 0.01915708812260536
$ws.Range("N15").Value = 0.007662835249042145
$ws.Range("O15").Value = 0.08812260536398467
$ws.Range("S15").Value = 0.2681992337164751
$ws.Range("F16").Value = 0.0223463687150838
$ws.Range("H16").Value = 0.1508379888268156
$ws.Range("I16").Value = 0.09497206703910614
$ws.Range("J16").Value = 0.3687150837988827
$ws.Range("K16").Value = 0.111731843575419
$ws.Range("M16").Value = 0.0223463687150838
$ws.Range("N16").Value = 0.0111731843575419
$ws.Range("O16").Value = 0.05027932960893855
$ws.Range("S16").Value = 0.1675977653631285
$ws.Range("F17").Value = 0.02669902912621359
$ws.Range("H17").Value = 0.2063106796116505
$ws.Range("I17").Value = 0.07524271844660194
$ws.Range("J17").Value = 0.3980582524271845
$ws.Range("K17").Value = 0.09466019417475728
$ws.Range("M17").Value = 0.009708737864077669
$ws.Range("O17").Value = 0.04611650485436893
$ws.Range("S17").Value = 0.1432038834951456
$ws.Range("F18").Value = 0.02973977695167286
$ws.Range("H18").Value = 0.1784386617100372
$ws.Range("I18").Value = 0.09293680297397769
$ws.Range("J18").Value = 0.3940520446096654
$ws.Range("K18").Value = 0.08921933085501858
$ws.Range("M18").Value = 0.02602230483271376
$ws.Range("O18").Value = 0.05576208178438662
$ws.Range("S18").Value = 0.1338289962825279
$ws.Range("F19").Value = 0.01327433628318584
$ws.Range("H19").Value = 0.226401179941003
$ws.Range("I19").Value = 0.06563421828908554
$ws.Range("J19").Value = 0.3488200589970502
$ws.Range("K19").Value = 0.1216814159292035
$ws.Range("M19").Value = 0.02433628318584071
$ws.Range("N19").Value = 0.0007374631268436578
$ws.Range("O19").Value = 0.08112094395280237
$ws.Range("S19").Value = 0.1179941002949852
